$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 3 for columns A, Q, R
$aTemp = $ws.Range("A2").Value2
$ws.Range("A2").Value2 = $ws.Range("A3").Value2
$ws.Range("A3").Value2 = $aTemp

$qTemp = $ws.Range("Q2").Value2
$ws.Range("Q2").Value2 = $ws.Range("Q3").Value2
$ws.Range("Q3").Value2 = $qTemp

$rTemp = $ws.Range("R2").Value2
$ws.Range("R2").Value2 = $ws.Range("R3").Value2
$ws.Range("R3").Value2 = $rTemp
